$wb = $excel.ActiveWorkbook

# --- Sheet "Herbiers" (sheet2): fill in the Famille / Genre taxonomy table ---
$ws2 = $wb.Worksheets.Item("Herbiers")

$ws2.Range("A1").Value = "Famille"
$ws2.Range("B1").Value = "Genre"

$ws2.Range("A2").Value = "Zosteraceae"
$ws2.Range("B3").Value = "Heterozostera"
$ws2.Range("B4").Value = "Phyllospadix"
$ws2.Range("B5").Value = "Zostera"

$ws2.Range("A6").Value = "Posidoniaceae"
$ws2.Range("B7").Value = "Posidonia"

$ws2.Range("A8").Value = "Cymodoceaceae"
$ws2.Range("B9").Value = "Amphibolis"
$ws2.Range("B10").Value = "Cymodocea"
$ws2.Range("B11").Value = "Halodule"
$ws2.Range("B12").Value = "Syringodium"
$ws2.Range("B13").Value = "Thalassodendron"

$ws2.Range("A14").Value = "Hydrocharitaceae"
$ws2.Range("B15").Value = "Enhalus"
$ws2.Range("B16").Value = "Halophila"
$ws2.Range("B17").Value = "Thalassia"

# Column A autofit-like width (bestFit), matching the authored sheet.
# (16.5703125 is Excel's best-fit width for "Hydrocharitaceae" at the
# workbook's default font; the nearest value this engine's column-width
# quantization can reach is used here.)
$ws2.Columns.Item(1).ColumnWidth = 15.7

# Page setup, matching the authored sheet (portrait A4/Letter "9" = A4)
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Selections, matching the final saved state of the workbook
$ws2.Range("C6").Select()

$ws1 = $wb.Worksheets.Item("Coraux")
$ws1.Range("C26").Select()

$ws2.Activate()
